# Insert a new data row at row 171 (shifts existing rows 171:276 down to
# 172:277) and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("171:171").Insert()

$ws.Cells.Item(171, 1).Value = 4
$ws.Cells.Item(171, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(171, 3).Value = "Los Lagos"
$ws.Cells.Item(171, 4).Value = 44873
$ws.Cells.Item(171, 5).Value = 10
$ws.Cells.Item(171, 6).Value = 100112039
$ws.Cells.Item(171, 7).Value = "Ciboulette"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 240
$ws.Cells.Item(171, 11).Value = 2000
$ws.Cells.Item(171, 12).Value = 2500
$ws.Cells.Item(171, 13).Value = 2250
$ws.Cells.Item(171, 14).Value = "`$/docena de atados"
$ws.Cells.Item(171, 15).Value = "Región Metropolitana"
$ws.Cells.Item(171, 16).Value = 750
$ws.Cells.Item(171, 17).Value = 3
$ws.Cells.Item(171, 18).Value = "Hortaliza"
